$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (SCD0282 -> SCD0018)
$ws.Name = "SCD0018"

# Update TC_ID column (B) from "DGS-297" to "SCD0018-005" on all three data rows
$ws.Range("B2").Value = "SCD0018-005"
$ws.Range("B3").Value = "SCD0018-005"
$ws.Range("B4").Value = "SCD0018-005"

# Widen column B to fit the new, longer TC_ID text
$ws.Columns.Item(2).ColumnWidth = 12

# Restore the view to the top-left and move the selection to B5
$ws.Range("B5").Select()
